# Apply the edit described by the diff: reshuffle the stimulus rows (2..27)
# for sheet '9_scenecat_categorization_living_rooms_2': column F (trial_total)
# decreases by 28 uniformly, and columns H, I, K, L, M, N, O, P, Q, R, S get new
# values per row (a row-level permutation of the trial data -- two old rows are
# dropped and two brand-new stimuli are introduced).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 27; $r++) {
    $oldTrialTotal = $ws.Cells.Item($r, 6).Value2   # column F: trial_total
    $ws.Cells.Item($r, 6).Value = $oldTrialTotal - 28
}

# Row 2
$ws.Cells.Item(2, 8).Value = "living_rooms"
$ws.Cells.Item(2, 9).Value = "target"
$ws.Cells.Item(2, 11).Value = "j"
$ws.Cells.Item(2, 12).Value = "stimuli/img_amsgw.png"
$ws.Cells.Item(2, 13).Value = 86.08510638297872
$ws.Cells.Item(2, 14).Value = 65.95744680851064
$ws.Cells.Item(2, 15).Value = 76.02127659574468
$ws.Cells.Item(2, 16).Value = 47
$ws.Cells.Item(2, 17).Value = 9
$ws.Cells.Item(2, 18).Value = 9
$ws.Cells.Item(2, 19).Value = 9

# Row 3
$ws.Cells.Item(3, 8).Value = "bedrooms"
$ws.Cells.Item(3, 9).Value = "distractor"
$ws.Cells.Item(3, 11).Value = "f"
$ws.Cells.Item(3, 12).Value = "stimuli/img_91csq.png"
$ws.Cells.Item(3, 13).Value = 50.44736842105263
$ws.Cells.Item(3, 14).Value = 28.34210526315789
$ws.Cells.Item(3, 15).Value = 39.39473684210526
$ws.Cells.Item(3, 16).Value = 38
$ws.Cells.Item(3, 17).Value = 2
$ws.Cells.Item(3, 18).Value = 2
$ws.Cells.Item(3, 19).Value = 2

# Row 4
$ws.Cells.Item(4, 8).Value = "living_rooms"
$ws.Cells.Item(4, 9).Value = "target"
$ws.Cells.Item(4, 11).Value = "j"
$ws.Cells.Item(4, 12).Value = "stimuli/img_9oofc.png"
$ws.Cells.Item(4, 13).Value = 82.47619047619048
$ws.Cells.Item(4, 14).Value = 65.5
$ws.Cells.Item(4, 15).Value = 73.98809523809524
$ws.Cells.Item(4, 16).Value = 42
$ws.Cells.Item(4, 17).Value = 8
$ws.Cells.Item(4, 18).Value = 8
$ws.Cells.Item(4, 19).Value = 8

# Row 5
$ws.Cells.Item(5, 8).Value = "living_rooms"
$ws.Cells.Item(5, 9).Value = "target"
$ws.Cells.Item(5, 11).Value = "j"
$ws.Cells.Item(5, 12).Value = "stimuli/img_rru0v.png"
$ws.Cells.Item(5, 13).Value = 56.45238095238095
$ws.Cells.Item(5, 14).Value = 39.42857142857143
$ws.Cells.Item(5, 15).Value = 47.94047619047619
$ws.Cells.Item(5, 16).Value = 42
$ws.Cells.Item(5, 17).Value = 4
$ws.Cells.Item(5, 18).Value = 4
$ws.Cells.Item(5, 19).Value = 4

# Row 6
$ws.Cells.Item(6, 8).Value = "living_rooms"
$ws.Cells.Item(6, 9).Value = "target"
$ws.Cells.Item(6, 11).Value = "j"
$ws.Cells.Item(6, 12).Value = "stimuli/img_2qhro.png"
$ws.Cells.Item(6, 13).Value = 81.73809523809524
$ws.Cells.Item(6, 14).Value = 62.73809523809524
$ws.Cells.Item(6, 15).Value = 72.23809523809524
$ws.Cells.Item(6, 16).Value = 42
$ws.Cells.Item(6, 17).Value = 8
$ws.Cells.Item(6, 18).Value = 8
$ws.Cells.Item(6, 19).Value = 8

# Row 7
$ws.Cells.Item(7, 8).Value = "bedrooms"
$ws.Cells.Item(7, 9).Value = "distractor"
$ws.Cells.Item(7, 11).Value = "f"
$ws.Cells.Item(7, 12).Value = "stimuli/img_d3t0o.png"
$ws.Cells.Item(7, 13).Value = 66.95121951219512
$ws.Cells.Item(7, 14).Value = 42.92682926829269
$ws.Cells.Item(7, 15).Value = 54.9390243902439
$ws.Cells.Item(7, 16).Value = 41
$ws.Cells.Item(7, 17).Value = 4
$ws.Cells.Item(7, 18).Value = 4
$ws.Cells.Item(7, 19).Value = 4

# Row 8
$ws.Cells.Item(8, 8).Value = "living_rooms"
$ws.Cells.Item(8, 9).Value = "target"
$ws.Cells.Item(8, 11).Value = "j"
$ws.Cells.Item(8, 12).Value = "stimuli/img_iudc4.png"
$ws.Cells.Item(8, 13).Value = 73.625
$ws.Cells.Item(8, 14).Value = 52.275
$ws.Cells.Item(8, 15).Value = 62.95
$ws.Cells.Item(8, 16).Value = 40
$ws.Cells.Item(8, 17).Value = 6
$ws.Cells.Item(8, 18).Value = 6
$ws.Cells.Item(8, 19).Value = 6

# Row 9
$ws.Cells.Item(9, 8).Value = "bedrooms"
$ws.Cells.Item(9, 9).Value = "distractor"
$ws.Cells.Item(9, 11).Value = "f"
$ws.Cells.Item(9, 12).Value = "stimuli/img_jp28n.png"
$ws.Cells.Item(9, 13).Value = 65.02564102564102
$ws.Cells.Item(9, 14).Value = 44.97435897435897
$ws.Cells.Item(9, 15).Value = 55
$ws.Cells.Item(9, 16).Value = 39
$ws.Cells.Item(9, 17).Value = 4
$ws.Cells.Item(9, 18).Value = 4
$ws.Cells.Item(9, 19).Value = 4

# Row 10
$ws.Cells.Item(10, 8).Value = "bedrooms"
$ws.Cells.Item(10, 9).Value = "distractor"
$ws.Cells.Item(10, 11).Value = "f"
$ws.Cells.Item(10, 12).Value = "stimuli/img_5il0t.png"
$ws.Cells.Item(10, 13).Value = 48.09523809523809
$ws.Cells.Item(10, 14).Value = 30.90476190476191
$ws.Cells.Item(10, 15).Value = 39.5
$ws.Cells.Item(10, 16).Value = 42
$ws.Cells.Item(10, 17).Value = 2
$ws.Cells.Item(10, 18).Value = 2
$ws.Cells.Item(10, 19).Value = 2

# Row 11
$ws.Cells.Item(11, 8).Value = "living_rooms"
$ws.Cells.Item(11, 9).Value = "target"
$ws.Cells.Item(11, 11).Value = "j"
$ws.Cells.Item(11, 12).Value = "stimuli/img_ac0ey.png"
$ws.Cells.Item(11, 13).Value = 86.62222222222222
$ws.Cells.Item(11, 14).Value = 70.02222222222223
$ws.Cells.Item(11, 15).Value = 78.32222222222222
$ws.Cells.Item(11, 16).Value = 45
$ws.Cells.Item(11, 17).Value = 9
$ws.Cells.Item(11, 18).Value = 9
$ws.Cells.Item(11, 19).Value = 9

# Row 12
$ws.Cells.Item(12, 8).Value = "living_rooms"
$ws.Cells.Item(12, 9).Value = "target"
$ws.Cells.Item(12, 11).Value = "j"
$ws.Cells.Item(12, 12).Value = "stimuli/img_xr3up.png"
$ws.Cells.Item(12, 13).Value = 76.24444444444444
$ws.Cells.Item(12, 14).Value = 55.88888888888889
$ws.Cells.Item(12, 15).Value = 66.06666666666666
$ws.Cells.Item(12, 16).Value = 45
$ws.Cells.Item(12, 17).Value = 7
$ws.Cells.Item(12, 18).Value = 7
$ws.Cells.Item(12, 19).Value = 7

# Row 13
$ws.Cells.Item(13, 8).Value = "living_rooms"
$ws.Cells.Item(13, 9).Value = "target"
$ws.Cells.Item(13, 11).Value = "j"
$ws.Cells.Item(13, 12).Value = "stimuli/img_zxvl3.png"
$ws.Cells.Item(13, 13).Value = 68.78260869565217
$ws.Cells.Item(13, 14).Value = 47.56521739130435
$ws.Cells.Item(13, 15).Value = 58.17391304347827
$ws.Cells.Item(13, 16).Value = 46
$ws.Cells.Item(13, 17).Value = 5
$ws.Cells.Item(13, 18).Value = 5
$ws.Cells.Item(13, 19).Value = 5

# Row 14
$ws.Cells.Item(14, 8).Value = "living_rooms"
$ws.Cells.Item(14, 9).Value = "target"
$ws.Cells.Item(14, 11).Value = "j"
$ws.Cells.Item(14, 12).Value = "stimuli/img_tn8ys.png"
$ws.Cells.Item(14, 13).Value = 86.70454545454545
$ws.Cells.Item(14, 14).Value = 72.40909090909091
$ws.Cells.Item(14, 15).Value = 79.55681818181819
$ws.Cells.Item(14, 16).Value = 44
$ws.Cells.Item(14, 17).Value = 10
$ws.Cells.Item(14, 18).Value = 10
$ws.Cells.Item(14, 19).Value = 10

# Row 15
$ws.Cells.Item(15, 8).Value = "living_rooms"
$ws.Cells.Item(15, 9).Value = "target"
$ws.Cells.Item(15, 11).Value = "j"
$ws.Cells.Item(15, 12).Value = "stimuli/img_kq9s9.png"
$ws.Cells.Item(15, 13).Value = 62.30232558139535
$ws.Cells.Item(15, 14).Value = 39.97674418604651
$ws.Cells.Item(15, 15).Value = 51.13953488372093
$ws.Cells.Item(15, 16).Value = 43
$ws.Cells.Item(15, 17).Value = 4
$ws.Cells.Item(15, 18).Value = 4
$ws.Cells.Item(15, 19).Value = 4

# Row 16
$ws.Cells.Item(16, 8).Value = "living_rooms"
$ws.Cells.Item(16, 9).Value = "target"
$ws.Cells.Item(16, 11).Value = "j"
$ws.Cells.Item(16, 12).Value = "stimuli/img_rg4in.png"
$ws.Cells.Item(16, 13).Value = 49.3695652173913
$ws.Cells.Item(16, 14).Value = 30.21739130434782
$ws.Cells.Item(16, 15).Value = 39.79347826086956
$ws.Cells.Item(16, 16).Value = 46
$ws.Cells.Item(16, 17).Value = 3
$ws.Cells.Item(16, 18).Value = 3
$ws.Cells.Item(16, 19).Value = 3

# Row 17
$ws.Cells.Item(17, 8).Value = "living_rooms"
$ws.Cells.Item(17, 9).Value = "target"
$ws.Cells.Item(17, 11).Value = "j"
$ws.Cells.Item(17, 12).Value = "stimuli/img_of8d6.png"
$ws.Cells.Item(17, 13).Value = 26.04878048780488
$ws.Cells.Item(17, 14).Value = 19.14634146341463
$ws.Cells.Item(17, 15).Value = 22.59756097560975
$ws.Cells.Item(17, 16).Value = 41
$ws.Cells.Item(17, 17).Value = 1
$ws.Cells.Item(17, 18).Value = 1
$ws.Cells.Item(17, 19).Value = 1

# Row 18
$ws.Cells.Item(18, 8).Value = "bedrooms"
$ws.Cells.Item(18, 9).Value = "distractor"
$ws.Cells.Item(18, 11).Value = "f"
$ws.Cells.Item(18, 12).Value = "stimuli/img_twj5p.png"
$ws.Cells.Item(18, 13).Value = 67.71739130434783
$ws.Cells.Item(18, 14).Value = 42.08695652173913
$ws.Cells.Item(18, 15).Value = 54.90217391304348
$ws.Cells.Item(18, 16).Value = 46
$ws.Cells.Item(18, 17).Value = 4
$ws.Cells.Item(18, 18).Value = 4
$ws.Cells.Item(18, 19).Value = 4

# Row 19
$ws.Cells.Item(19, 8).Value = "living_rooms"
$ws.Cells.Item(19, 9).Value = "target"
$ws.Cells.Item(19, 11).Value = "j"
$ws.Cells.Item(19, 12).Value = "stimuli/img_vh7v8.png"
$ws.Cells.Item(19, 13).Value = 78.70454545454545
$ws.Cells.Item(19, 14).Value = 59.63636363636363
$ws.Cells.Item(19, 15).Value = 69.17045454545455
$ws.Cells.Item(19, 16).Value = 44
$ws.Cells.Item(19, 17).Value = 7
$ws.Cells.Item(19, 18).Value = 7
$ws.Cells.Item(19, 19).Value = 7

# Row 20
$ws.Cells.Item(20, 8).Value = "living_rooms"
$ws.Cells.Item(20, 9).Value = "target"
$ws.Cells.Item(20, 11).Value = "j"
$ws.Cells.Item(20, 12).Value = "stimuli/img_ra2nm.png"
$ws.Cells.Item(20, 13).Value = 70.75
$ws.Cells.Item(20, 14).Value = 50.375
$ws.Cells.Item(20, 15).Value = 60.5625
$ws.Cells.Item(20, 16).Value = 40
$ws.Cells.Item(20, 17).Value = 6
$ws.Cells.Item(20, 18).Value = 6
$ws.Cells.Item(20, 19).Value = 6

# Row 21
$ws.Cells.Item(21, 8).Value = "living_rooms"
$ws.Cells.Item(21, 9).Value = "target"
$ws.Cells.Item(21, 11).Value = "j"
$ws.Cells.Item(21, 12).Value = "stimuli/img_vgh2g.png"
$ws.Cells.Item(21, 13).Value = 93.81395348837209
$ws.Cells.Item(21, 14).Value = 78.27906976744185
$ws.Cells.Item(21, 15).Value = 86.04651162790697
$ws.Cells.Item(21, 16).Value = 43
$ws.Cells.Item(21, 17).Value = 10
$ws.Cells.Item(21, 18).Value = 10
$ws.Cells.Item(21, 19).Value = 10

# Row 22
$ws.Cells.Item(22, 8).Value = "living_rooms"
$ws.Cells.Item(22, 9).Value = "target"
$ws.Cells.Item(22, 11).Value = "j"
$ws.Cells.Item(22, 12).Value = "stimuli/img_j4ttn.png"
$ws.Cells.Item(22, 13).Value = 12.61904761904762
$ws.Cells.Item(22, 14).Value = 11.42857142857143
$ws.Cells.Item(22, 15).Value = 12.02380952380952
$ws.Cells.Item(22, 16).Value = 42
$ws.Cells.Item(22, 17).Value = 1
$ws.Cells.Item(22, 18).Value = 1
$ws.Cells.Item(22, 19).Value = 1

# Row 23
$ws.Cells.Item(23, 8).Value = "living_rooms"
$ws.Cells.Item(23, 9).Value = "target"
$ws.Cells.Item(23, 11).Value = "j"
$ws.Cells.Item(23, 12).Value = "stimuli/img_swq34.png"
$ws.Cells.Item(23, 13).Value = 64.11363636363636
$ws.Cells.Item(23, 14).Value = 43.04545454545455
$ws.Cells.Item(23, 15).Value = 53.57954545454545
$ws.Cells.Item(23, 16).Value = 44
$ws.Cells.Item(23, 17).Value = 5
$ws.Cells.Item(23, 18).Value = 5
$ws.Cells.Item(23, 19).Value = 5

# Row 24
$ws.Cells.Item(24, 8).Value = "living_rooms"
$ws.Cells.Item(24, 9).Value = "target"
$ws.Cells.Item(24, 11).Value = "j"
$ws.Cells.Item(24, 12).Value = "stimuli/img_rych7.png"
$ws.Cells.Item(24, 13).Value = 30.4468085106383
$ws.Cells.Item(24, 14).Value = 23.4468085106383
$ws.Cells.Item(24, 15).Value = 26.9468085106383
$ws.Cells.Item(24, 16).Value = 47
$ws.Cells.Item(24, 17).Value = 2
$ws.Cells.Item(24, 18).Value = 2
$ws.Cells.Item(24, 19).Value = 2

# Row 25
$ws.Cells.Item(25, 8).Value = "living_rooms"
$ws.Cells.Item(25, 9).Value = "target"
$ws.Cells.Item(25, 11).Value = "j"
$ws.Cells.Item(25, 12).Value = "stimuli/img_24rt2.png"
$ws.Cells.Item(25, 13).Value = 55.26829268292683
$ws.Cells.Item(25, 14).Value = 34.19512195121951
$ws.Cells.Item(25, 15).Value = 44.73170731707317
$ws.Cells.Item(25, 16).Value = 41
$ws.Cells.Item(25, 17).Value = 3
$ws.Cells.Item(25, 18).Value = 3
$ws.Cells.Item(25, 19).Value = 3

# Row 26
$ws.Cells.Item(26, 8).Value = "bedrooms"
$ws.Cells.Item(26, 9).Value = "distractor"
$ws.Cells.Item(26, 11).Value = "f"
$ws.Cells.Item(26, 12).Value = "stimuli/img_3h4c9.png"
$ws.Cells.Item(26, 13).Value = 85.47619047619048
$ws.Cells.Item(26, 14).Value = 67.26190476190476
$ws.Cells.Item(26, 15).Value = 76.36904761904762
$ws.Cells.Item(26, 16).Value = 42
$ws.Cells.Item(26, 17).Value = 9
$ws.Cells.Item(26, 18).Value = 9
$ws.Cells.Item(26, 19).Value = 9

# Row 27
$ws.Cells.Item(27, 8).Value = "living_rooms"
$ws.Cells.Item(27, 9).Value = "target"
$ws.Cells.Item(27, 11).Value = "j"
$ws.Cells.Item(27, 12).Value = "stimuli/img_syam3.png"
$ws.Cells.Item(27, 13).Value = 41.32432432432432
$ws.Cells.Item(27, 14).Value = 26.2972972972973
$ws.Cells.Item(27, 15).Value = 33.81081081081081
$ws.Cells.Item(27, 16).Value = 37
$ws.Cells.Item(27, 17).Value = 2
$ws.Cells.Item(27, 18).Value = 2
$ws.Cells.Item(27, 19).Value = 2
